$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Beer")
$ws.Activate()

# Delete entire column E (SetPointHigh) - shifts remaining columns left
$ws.Columns("E:E").Delete()

# Delete rows 22 and 23 (the two TIC-110 / TIC-109 rows) - shifts remaining rows up
$ws.Rows("22:23").Delete()

# Leave the freshly-shifted rows 22:23 selected, matching the post-edit state
$ws.Range("A22:XFD23").Select()
